# Apply updated cryptocurrency price/volume data to Sheet1
# Values are written as text (quote-prefixed) to preserve formats such as
# "42.721.78" or "1.00" that Excel would otherwise coerce into numbers,
# then the cell's Style is reset to 'Normal' so no stray formatting is kept.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = "'"  # single quote used as a text-prefix so Excel stores the value as a string

$ws.Range('D2').Value = $q + '42.721.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = $q + '  -0.91%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = $q + '2.537.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = $q + '  -1.06%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = $q + '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = $q + '  +0.07%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = $q + '308.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = $q + '  -2.17%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = $q + '100.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = $q + '  +4.05%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = $q + '  -1.17%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D9').Value = $q + '0.529'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = $q + '  -2.11%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = $q + '36.03'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = $q + '  +1.50%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = $q + '0.0805'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = $q + '  -1.01%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = $q + '7.35'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = $q + '  -1.15%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = $q + '  +0.01%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = $q + '2.930.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = $q + '  -0.91%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = $q + '  +5.63%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = $q + '2.604.06'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = $q + '  +0.36%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = $q + '  -3.33%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = $q + '42.709.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = $q + '  -0.94%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E20').Value = $q + '  -0.79%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = $q + '12.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = $q + '  -2.59%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = $q + '69.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = $q + '  +0.41%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = $q + '243.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = $q + '  -3.78%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = $q + '2.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = $q + '  -1.98%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = $q + '  -1.36%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = $q + '  +0.01%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = $q + '26.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = $q + '  -2.84%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = $q + '2.35'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = $q + '  -3.64%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = $q + '39.31'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = $q + '  -1.92%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = $q + '10.15'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = $q + '  -0.25%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = $q + 'Filecoin'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = $q + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = $q + '5.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = $q + '  -0.65%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = $q + 'Monero'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = $q + 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = $q + '156.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = $q + '  +0.71%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = $q + '2.76'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = $q + '  +12.67%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = $q + '0.0793'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = $q + '  -1.68%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = $q + '  -2.80%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = $q + '18.36'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = $q + '  -3.09%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = $q + '  -4.26%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = $q + '  -6.36%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = $q + '  +0.01%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = $q + '  +0.63%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = $q + '4.32'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = $q + '  +8.77%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = $q + '21.96'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = $q + '  -2.09%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = $q + '  +0.07%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = $q + '  +1.90%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = $q + '0.0298'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = $q + '  -1.91%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = $q + '1.972.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = $q + '  -1.52%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = $q + '  -0.56%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = $q + '81.22'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = $q + '  -1.68%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = $q + 'Algorand'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = $q + 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = $q + '0.192'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = $q + '  -0.50%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = $q + 'SEI'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = $q + 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = $q + '0.858'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = $q + '  +10.56%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = $q + '2.725.36'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = $q + '  -2.95%  '
$ws.Range('E51').Style = 'Normal'
